# Daily GitHub Actions refresh of the cryptos price/volume table.
# Price (D) and Volume(1h) (E) columns are plain text in this sheet (not
# numbers), and two coin pairs (rows 15/16 and 47/48) swapped rank order.
# A leading apostrophe forces Excel to keep "clean" numeric-looking prices
# (e.g. "0.999", "27.10", "1.00") stored as literal text instead of being
# auto-coerced into a Double (which would silently drop trailing zeros /
# the decimal point).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.501.77"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "3.557.34"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'598.56"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "'141.08"
$ws.Range("E6").Value = "  +4.28%  "
$ws.Range("D7").Value = "3.557.18"
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'0.127"
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("D11").Value = "'7.19"
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = "  +3.96%  "
$ws.Range("D13").Value = "4.153.45"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").Value = "'0.0000189"
$ws.Range("E14").Value = "  +4.83%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.563.58"
$ws.Range("E15").Value = "  +4.26%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'27.10"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "65.341.34"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'10.36"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("D20").Value = "'5.88"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").Value = "'14.25"
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("D22").Value = "'397.54"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").Value = "'0.572"
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("D24").Value = "'74.63"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "3.694.58"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +9.22%  "
$ws.Range("D28").Value = "'7.82"
$ws.Range("E28").Value = "  +7.33%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "3.565.50"
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("D33").Value = "'24.10"
$ws.Range("E33").Value = "  +6.23%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'7.07"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("D38").Value = "'168.66"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "'4.99"
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("D41").Value = "'0.0809"
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("D42").Value = "'0.828"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'26.55"
$ws.Range("E43").Value = "  +16.76%  "
$ws.Range("D44").Value = "'42.66"
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("D46").Value = "'4.46"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.68"
$ws.Range("E47").Value = "  +3.93%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.19"
$ws.Range("E48").Value = "  +8.37%  "
$ws.Range("D49").Value = "'6.84"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("D50").Value = "2.396.41"
$ws.Range("E50").Value = "  +9.56%  "
$ws.Range("D51").Value = "'2.14"
$ws.Range("E51").Value = "  +0.62%  "
